# Update header labels to English (with NER-friendly naming)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Entry"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Filename"

# Widen column B to fit the longer "Name" header content
$ws.Columns.Item(2).ColumnWidth = 36.5

# Move/collapse the selection to I17
$ws.Range("I17").Select()
